$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 62: new date entry (2023/1/10) ---
# Copy formatting from an existing date cell (A1, style s="1" -> m/d/yyyy)
$ws.Range("A1").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A62").Value = 44936

# --- Row 63: time entry + "HW示範" note ---
# Copy formatting from an existing time cell (A3, style s="2" -> h:mm)
$ws.Range("A3").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A63").Value = 0.44791666666666669
$ws.Range("B63").Value = "HW示範"

# --- Row 64: time entry + "var 匿名物件 介面" note ---
$ws.Range("A3").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A64").Value = 0.4861111111111111
$ws.Range("B64").Value = "var 匿名物件 介面"

$excel.CutCopyMode = $false

# --- Column A width (book author widened col A / best-fit) ---
$ws.Range("A1:A64").ColumnWidth = 8.71

# --- Update view: scrolled down, new selection at A65 ---
$ws.Range("A65").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 59
$win.ScrollColumn = 1
